$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CTC" (sheet2.xml): change prescale (A4) and timer-bit divisor (A8)
# ---------------------------------------------------------------------------
$ctc = $wb.Worksheets.Item("CTC")
$ctc.Range("A4").Value = 0.5
$ctc.Range("A8").Value = 256
$ctc.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet "RPM in CTC" (sheet3.xml)
# ---------------------------------------------------------------------------
$rpm = $wb.Worksheets.Item("RPM in CTC")

# Rename the label so the "Interrupts per Second" row becomes "... [Hz]"
$rpm.Range("B7").Value = "Interrupts per Second [Hz]"

# Update inputs
$rpm.Range("A4").Value = 1
$rpm.Range("A6").Value = 16
$rpm.Range("A10").Value = 256

# Change the interrupts-per-minute formula from *60 to /60
$rpm.Range("B8").Formula = "=B6/60"

# New row: total interrupt frequency in Hz
$rpm.Range("D12").Formula = "=B10*A8*1000000"
$rpm.Range("D12").NumberFormat = "0.00E+00"

# New (empty) cell carrying the custom "0.0E+00" scientific format
$rpm.Range("D14").NumberFormat = "0.0E+00"

# B10 (the compare-match time) keeps a scientific format, now the plain
# built-in 0.00E+00 instead of the old custom 0.000E+00
$rpm.Range("B10").NumberFormat = "0.00E+00"

# Column B needs to be a touch wider to fit the new label text
$rpm.Columns.Item(2).ColumnWidth = 24.42578125

$rpm.Range("C8").Select()
$rpm.Activate()
